# Adds four new observation rows (7-10) to the "Artfynd" sheet, mirroring
# the existing rows' column layout (same species "Knärot" / Goodyera repens
# record group, reported by Kim Hultgren on 2023-09-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        A = 112092130; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        P = "Sollefteå (Sollefteå), Ång"
        Q = 584352.4882331375; R = 7048231.676015709; S = 25
        T = "Västernorrland"; U = "Sollefteå"; V = "Ångermanland"; W = "Ramsele"
        Y = "2023-09-14"; Z = "17:22"; AA = "2023-09-14"; AB = "17:22"
        AD = $false; AE = $false; AG = $false
        AW = "Kim Hultgren"; AX = "Kim Hultgren"
    },
    @{
        A = 112092161; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        I = "10"
        P = "Sollefteå (Sollefteå), Ång"
        Q = 584329.919621415; R = 7048274.339291978; S = 25
        T = "Västernorrland"; U = "Sollefteå"; V = "Ångermanland"; W = "Ramsele"
        Y = "2023-09-14"; Z = "17:22"; AA = "2023-09-14"; AB = "17:22"
        AD = $false; AE = $false; AG = $false
        AW = "Kim Hultgren"; AX = "Kim Hultgren"
    },
    @{
        A = 112092586; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        I = "20"
        P = "Sollefteå (Sollefteå), Ång"
        Q = 584400.9675979441; R = 7048356.949537945; S = 25
        T = "Västernorrland"; U = "Sollefteå"; V = "Ångermanland"; W = "Ramsele"
        Y = "2023-09-14"; Z = "17:46"; AA = "2023-09-14"; AB = "17:46"
        AD = $false; AE = $false; AG = $false
        AW = "Kim Hultgren"; AX = "Kim Hultgren"
    },
    @{
        A = 112092066; B = 96348; C = "Ovaliderad"; D = "VU"; E = 220787
        F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."
        P = "Sollefteå (Sollefteå), Ång"
        Q = 584345.5636095351; R = 7048206.515963284; S = 25
        T = "Västernorrland"; U = "Sollefteå"; V = "Ångermanland"; W = "Ramsele"
        Y = "2023-09-14"; Z = "17:18"; AA = "2023-09-14"; AB = "17:18"
        AD = $false; AE = $false; AG = $false
        AW = "Kim Hultgren"; AX = "Kim Hultgren"
    }
)

# Column letter -> 1-based column index, in the same order the source
# worksheet populates each record row.
$textCols = @{
    C = 3; D = 4; F = 6; G = 7; H = 8; I = 9
    P = 16; T = 20; U = 21; V = 22; W = 23
    Y = 25; Z = 26; AA = 27; AB = 28
    AW = 49; AX = 50
}
$numCols = @{ A = 1; B = 2; E = 5; Q = 17; R = 18; S = 19 }
$boolCols = @{ AD = 30; AE = 31; AG = 33 }

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $data = $newRows[$i]

    foreach ($col in $numCols.Keys) {
        if ($data.ContainsKey($col)) {
            $ws.Cells.Item($rowNum, $numCols[$col]).Value = $data[$col]
        }
    }

    foreach ($col in $boolCols.Keys) {
        if ($data.ContainsKey($col)) {
            $ws.Cells.Item($rowNum, $boolCols[$col]).Value = $data[$col]
        }
    }

    # Force text storage for string columns so values that look like
    # dates/numbers (e.g. "2023-09-14", "17:22", "10") are not
    # auto-converted by Excel's type inference.
    foreach ($col in $textCols.Keys) {
        if ($data.ContainsKey($col)) {
            $cell = $ws.Cells.Item($rowNum, $textCols[$col])
            $cell.NumberFormat = "@"
            $cell.Value = $data[$col]
        }
    }
}
